$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.687.26"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "2.519.42"
$ws.Range("E3").Value = "  +9.19%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'298.07"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "'96.52"
$ws.Range("E6").Value = "  -3.32%  "

$ws.Range("D7").Value = "'0.584"
$ws.Range("E7").Value = "  +1.86%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  +5.83%  "

$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  +3.63%  "

$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").Value = "'7.55"
$ws.Range("E12").Value = "  +6.45%  "

$ws.Range("D13").Value = "2.905.79"
$ws.Range("E13").Value = "  +10.76%  "

$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "2.515.97"
$ws.Range("E15").Value = "  +9.25%  "

$ws.Range("D16").Value = "'0.879"
$ws.Range("E16").Value = "  +9.64%  "

$ws.Range("D17").Value = "'14.50"
$ws.Range("E17").Value = "  +5.57%  "

$ws.Range("D18").Value = "45.797.89"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "'13.20"
$ws.Range("E19").Value = "  +4.55%  "

$ws.Range("D20").Value = "0.0₃0962"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("D21").Value = "'6.49"
$ws.Range("E21").Value = "  +10.40%  "

$ws.Range("D22").Value = "'68.46"
$ws.Range("E22").Value = "  +3.64%  "

$ws.Range("D23").Value = "'249.48"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").Value = "'2.84"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("E25").Value = "  +8.09%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").Value = "'40.54"
$ws.Range("E27").Value = "  -1.63%  "

$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").Value = "'22.59"
$ws.Range("E30").Value = "  +11.59%  "

$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  +13.06%  "

$ws.Range("D32").Value = "'5.74"
$ws.Range("E32").Value = "  +6.65%  "

$ws.Range("D33").Value = "'2.21"
$ws.Range("E33").Value = "  +30.82%  "

$ws.Range("D34").Value = "'2.77"
$ws.Range("E34").Value = "  -0.99%  "

$ws.Range("D35").Value = "'148.95"
$ws.Range("E35").Value = "  +2.35%  "

$ws.Range("D36").Value = "'0.0795"
$ws.Range("E36").Value = "  +2.36%  "

$ws.Range("E37").Value = "  +3.24%  "

$ws.Range("E38").Value = "  +1.24%  "

$ws.Range("D39").Value = "'15.68"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").Value = "'4.08"
$ws.Range("E40").Value = "  +5.57%  "

$ws.Range("D41").Value = "'0.0308"
$ws.Range("E41").Value = "  +3.28%  "

$ws.Range("D42").Value = "'3.35"
$ws.Range("E42").Value = "  +6.49%  "

$ws.Range("D43").Value = "2.021.54"
$ws.Range("E43").Value = "  +12.11%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").Value = "'91.76"
$ws.Range("E45").Value = "  -2.38%  "

$ws.Range("D46").Value = "'16.63"
$ws.Range("E46").Value = "  +33.99%  "

$ws.Range("D47").Value = "'1.80"
$ws.Range("E47").Value = "  -4.06%  "

$ws.Range("E48").Value = "  +11.62%  "

$ws.Range("D49").Value = "'8.78"
$ws.Range("E49").Value = "  +11.60%  "

$ws.Range("D50").Value = "2.755.46"
$ws.Range("E50").Value = "  +9.88%  "

$ws.Range("E51").Value = "  +3.24%  "
